# Revert "Add concentration to antibodies"
#
# Undo the prior change that introduced a concentration_value /
# concentration_unit pair of columns (plus their supporting
# "concentration_unit list" sheet and data validation) into the
# "Export as TSV" sheet, and bumped the schema version from 1 back to 2.
# This script puts things back the way they were: version 1, no
# concentration columns.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$dataSheet = $wb.Worksheets.Item("Export as TSV")

# Column H currently holds concentration_value, I holds concentration_unit,
# J holds conjugated_cat_number, and K holds conjugated_tag. Grab the
# comment text belonging to J1/K1 now, before those columns get shifted left
# by the EntireColumn.Delete() below (deleting columns re-indexes cell
# values automatically, but leaves cell comments anchored to their original
# addresses, so we fix the comments up by hand afterwards).
$conjugatedCatNumberComment = $dataSheet.Range("J1").Comment.Text()
$conjugatedTagComment = $dataSheet.Range("K1").Comment.Text()

# Delete the concentration_value / concentration_unit columns outright. This
# removes the header cells and their data validation rules, and shifts
# conjugated_cat_number/conjugated_tag left into H:I.
$dataSheet.Range("H1:I1").EntireColumn.Delete()

# The comment objects physically anchored at H1/I1 still exist (still
# carrying their old, now-stale text) -- overwrite them with the text that
# used to belong to J1/K1 (which now visually shifted into H1/I1).
$dataSheet.Range("H1").Comment.Text($conjugatedCatNumberComment) | Out-Null
$dataSheet.Range("I1").Comment.Text($conjugatedTagComment) | Out-Null

# J1/K1's comments are now orphaned (no cell data lives under them anymore);
# drop them.
$dataSheet.Range("J1").Comment.Delete() | Out-Null
$dataSheet.Range("K1").Comment.Delete() | Out-Null

# Revert the schema version advertised in the "version list" sheet (and
# referenced by the dropdown validation on column A) from 2 back to 1.
$versionSheet = $wb.Worksheets.Item("version list")
$versionCell = $versionSheet.Range("A1")
$versionCell.NumberFormat = "@"
$versionCell.Value = "1"
$versionCell.Style = "Normal"

# The list validation on column A caches a rendered error message that
# names the allowed value; refresh it to match the reverted version.
$versionValidation = $dataSheet.Range("A2:A1048576").Validation
$versionValidation.ErrorMessage = "Value must be one of: 1."

# Drop the now-unused "concentration_unit list" sheet (and its lookup value
# "ug/ml") entirely.
$wb.Worksheets.Item("concentration_unit list").Delete() | Out-Null

# Leave the originally active sheet selected.
$dataSheet.Activate()
